$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.295792818069458
$ws.Range("B1").Value = 2.478166341781616
$ws.Range("C1").Value = 2.515751600265503
$ws.Range("D1").Value = 3.217341661453247
$ws.Range("E1").Value = 2.404553651809692
